$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Append new journal rows (58-65) at the bottom of the existing table,
#    copying formatting from analogous existing rows so no stray styles are
#    introduced.
# ---------------------------------------------------------------------------

# Row 58 is a blank separator row, like rows 9 / 23 / 34 / 41 / 48.
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A58:G58").PasteSpecial(-4122) | Out-Null

# Rows 59-65 are regular (non date-header) data rows, like row 6.
$ws.Range("A6:G6").Copy() | Out-Null
$ws.Range("A59:G65").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the values for the new rows.
# ---------------------------------------------------------------------------

# Row 59
$ws.Range("B59").Value2 = 0.35416666666666669
$ws.Range("C59").Value2 = 0.39583333333333331
$ws.Range("D59").Formula = "=C59-B59"
$ws.Range("E59").Value2 = "Implémentation"
$ws.Range("F59").Value2 = "Màj de la db pour faciliter l'import des données"
$ws.Range("G59").Value2 = "mld à recréer / utilisation de fichier backup au lieu de .sql pour créer la db / type de colonne changé pour celestial_object.right_ascension"
$ws.Range("A59:G59").EntireRow.RowHeight = 30

# Row 60
$ws.Range("B60").Value2 = 0.40972222222222227
$ws.Range("C60").Value2 = 0.48958333333333331
$ws.Range("D60").Formula = "=C60-B60"
$ws.Range("E60").Value2 = "Implémentation"
$ws.Range("G60").Value2 = "fichier csv erroné: m45 n'a pas de n° ngc --> son nom commun a été pris --> modif. csv ??"
$ws.Range("F60").Value2 = "Création du script permettant de remplir la table celestiabl_object"

# Row 61
$ws.Range("B61").Value2 = 0.53125
$ws.Range("C61").Value2 = 0.55208333333333337
$ws.Range("D61").Formula = "=C61-B61"
$ws.Range("E61").Value2 = "Implémentation"
$ws.Range("F61").Value2 = "Téléchargement d'images pour la bdd"
$ws.Range("G61").Value2 = "httrack est trop long --> extension firefox DownThemAll lien: http://www.httrack.com/page/2/fr/index.html"
$ws.Range("A61:G61").EntireRow.RowHeight = 30

# Row 62
$ws.Range("B62").Value2 = 0.55208333333333337
$ws.Range("C62").Value2 = 0.5625
$ws.Range("D62").Formula = "=C62-B62"
$ws.Range("E62").Value2 = "Implémentation"
$ws.Range("F62").Value2 = "Téléchargement d'images pour la bdd"
$ws.Range("G62").Value2 = "Extension Firefox: DownThemAll -> réussite ! Choisir jpg ET png / suite: renommer les images utiles pour le site"
$ws.Range("A62:G62").EntireRow.RowHeight = 30

# Row 63
$ws.Range("B63").Value2 = 0.5625
$ws.Range("C63").Value2 = 0.62847222222222221
$ws.Range("D63").Formula = "=C63-B63"
$ws.Range("E63").Value2 = "Implémentation"
$ws.Range("F63").Value2 = "Création d'un script pour rename le images"
$ws.Range("G63").Value2 = "src: https://stackoverflow.com/questions/40904836/how-to-get-n-files-in-a-directory-order-by-last-modified-date -- https://stackoverflow.com/questions/52152228/how-to-filter-a-list-with-a-list-of-strings-in-powershell -- https://stackoverflow.com/questions/11816218/renaming-files-in-powershell-using-the-folder-name -- https://stackoverflow.com/questions/51818485/increment-variable-in-powershell-from-within-if-statement-within-a-foreach-loop"
$ws.Range("A63:G63").EntireRow.RowHeight = 90

# Row 64
$ws.Range("B64").Value2 = 0.63888888888888895
$ws.Range("C64").Value2 = 0.67708333333333337
$ws.Range("D64").Formula = "=C64-B64"
$ws.Range("E64").Value2 = "Implémentation"
$ws.Range("F64").Value2 = "Création d'un script pour rename le images"
$ws.Range("G64").Value2 = "Script créé mais impossible de renommer les images --> à terminer"

# Row 65
$ws.Range("B65").Value2 = 0.67708333333333337
$ws.Range("C65").Value2 = 0.70486111111111116
$ws.Range("D65").Formula = "=C65-B65"
$ws.Range("E65").Value2 = "Documentation"
$ws.Range("G65").Value2 = "Questions à poser à la cdp dans fichiers /analyse/questions.txt"
$ws.Range("F65").Value2 = "Mise à jour du dossier de projet"

# ---------------------------------------------------------------------------
# 3. Update the sheet view to reflect the scrolled / selected state that
#    results from having added the new rows at the bottom.
# ---------------------------------------------------------------------------

$ws.Range("A53").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("G65").Select() | Out-Null
